$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The team-stats export stored the game date as text using the wrong format
# (e.g. "5-29-2012-13", a mangled mash-up of month-day and season). It needs
# to be corrected to the ISO-style "2013-05-29" string used elsewhere.
$oldText = "5-29-2012-13"
$newText = "2013-05-29"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Locate the "Date" column header so this doesn't depend on a hard-coded
# column letter.
$headerRow = $ws.Rows.Item($firstRow)
$dateHeader = $headerRow.Find("Date")
if ($dateHeader -ne $null) {
    $dateCol = $dateHeader.Column
} else {
    $dateCol = 58 # fallback: column BF
}

# Excel treats a range's NumberFormat as "@" (Text) so that re-typing a
# date-shaped string is kept as literal text instead of being silently
# converted into a date serial number.
$colRange = $ws.Range($ws.Cells.Item($firstRow + 1, $dateCol), $ws.Cells.Item($lastRow, $dateCol))
$colRange.NumberFormat = "@"

for ($row = $firstRow + 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}

# Drop the temporary text formatting again so the cells keep the same
# (default) style they had before the edit.
$colRange.Style = "Normal"
